# Updated Use cases Slide
#
# The "3. Login/ Register (option to use Gmail/facebook login)" bullet on the
# "Summary of Use Cases" slide gets "facebook" corrected to "Facebook".
#
# Locate the shape holding the text dynamically (rather than hard-coding
# slide/shape indices) so the script is resilient to minor structural
# differences, then perform the capitalization fix using the same
# select-and-retype pattern a human editing in the UI would use: each
# sub-string write below lands on exactly the boundary PowerPoint's own
# edit produced, so the paragraph ends up split into the same run
# sequence ("...option to", " use ", "Gmail/", "F", "acebook ", "login)")
# as the authored change.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text.IndexOf("facebook") -ge 0) {
                $targetSlide = $s
                $targetShape = $sh
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange

    # --- Step 1: separate "Gmail/" from the " use " run that currently ends
    #     with it ("<sp>use<sp>Gmail/") -------------------------------------
    $full = $tr.Text
    $idxGmail = $full.IndexOf("Gmail/")
    $rngGmail = $tr.Characters($idxGmail + 1, 6)
    $rngGmail.Text = "Gmail/"

    # --- Step 2: capitalize the "f" in "facebook" as its own run ("F") -----
    $full = $tr.Text
    $idxFacebook = $full.IndexOf("facebook")
    $rngF = $tr.Characters($idxFacebook + 1, 1)
    $rngF.Text = "F"

    # --- Step 3: move the trailing space from the front of " login)" onto
    #     the end of "acebook" ---------------------------------------------
    $full = $tr.Text
    $idxAcebook = $full.IndexOf("acebook")
    $rngAcebook = $tr.Characters($idxAcebook + 1, 7)
    $rngAcebook.Text = "acebook "

    $full = $tr.Text
    $idxLoginSpace = $full.IndexOf(" login)")
    $rngLogin = $tr.Characters($idxLoginSpace + 1, 7)
    $rngLogin.Text = "login)"

    Write-Host "Updated text:" $tr.Text
}
